$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J2: "001" -> "002" (must stay text, not become the number 2) ---
# Forcing a numeric-looking string to stay text normally stamps the cell
# with a "Text" number format / quote-prefix style. Do that on a scratch
# cell instead, copy just the resulting value over to J2, then wipe the
# scratch cell so the only lasting change is J2's content.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "002"
$scratch.Copy()
$ws.Range("J2").PasteSpecial(-4163)
$scratch.Clear()
# K2 already holds "001" and is left untouched.

# --- N2: report date ---
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# --- O2:AD2 numeric figures ---
$ws.Range("O2").Value = -81292482.95
$ws.Range("P2").Value = -965.0446451346
$ws.Range("Q2").Value = 683298744.51
$ws.Range("R2").Value = 8111.6207856779
$ws.Range("S2").Value = 681583932.78
$ws.Range("T2").Value = 8091.2638004144
$ws.Range("U2").Value = 16595658.63
$ws.Range("V2").Value = 197.0114691074
$ws.Range("W2").Value = 677911.42
$ws.Range("X2").Value = 8.047666426299999
$ws.Range("Y2").Value = 14182252.79
$ws.Range("Z2").Value = 168.3612877141
$ws.Range("AA2").Value = 56256538.62
$ws.Range("AB2").Value = 667.8363039106
$ws.Range("AC2").Value = -8423701.779999999
$ws.Range("AD2").Value = -122.2929552369
